$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "IP Address 4" column (column L) with header + 5 data rows
$ws.Range("L2").Value = "IP Address 4"
$ws.Range("L3").Value = "13.201.71.98"
$ws.Range("L4").Value = "13.233.73.254"
$ws.Range("L5").Value = "13.233.93.60"
$ws.Range("L6").Value = "13.126.114.210"
$ws.Range("L7").Value = "52.66.196.55"

# Match the header style used by the other header cells in row 2
$ws.Range("L2").Font.Bold = $true
$ws.Range("L2").WrapText = $true
$ws.Range("L2").VerticalAlignment = -4108

# Adjust column widths: K narrower, new L column width set
$ws.Columns.Item(11).ColumnWidth = 16.77734375
$ws.Columns.Item(12).ColumnWidth = 22.88671875

# Update the active selection to match the saved workbook state
$ws.Range("H11").Select()
